$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of scanned codes ("triangulos / circulos" nesting data) appended
# below the existing table. Force the range to Text format first so that
# numeric-looking codes (e.g. "123456", "1585") and the DD/MM/YYYY HH:MM:SS
# timestamps are stored as literal text, matching the rest of the sheet.
$newRange = $ws.Range("A28:C37")
$newRange.NumberFormat = "@"

$data = @(
    @("VDS1832", "15/10/2025 07:28:26", "123456"),
    @("VDS1833", "15/10/2025 07:48:59", "1922"),
    @("VDS1834", "15/10/2025 10:09:49", "Chorão-as-rosas"),
    @("VDS1835", "15/10/2025 10:11:46", "Noroaco"),
    @("VDS1836", "15/10/2025 10:15:13", "Noroaco"),
    @("VDS1837", "15/10/2025 10:19:27", "1585"),
    @("VDS1838", "15/10/2025 10:19:39", "1585"),
    @("VDS1839", "15/10/2025 10:19:53", "1585"),
    @("VDS1840", "15/10/2025 10:42:19", "1585"),
    @("VDS1841", "15/10/2025 10:42:43", "1585")
)

$row = 28
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
